$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$f2 = @'
<rpc-reply message-id="urn:uuid:ff35cb2e-873a-4fbc-bb04-d39cbbcc0c6d">
  <data>
    <network-instances>
      <network-instance>
        <name>Prueba_LxVPN</name>
        <config>
          <name>Prueba_LxVPN</name>
          <type>oc-ni-types:L3VRF</type>
        </config>
        <interfaces>
          <interface>
            <id>GigabitEthernet0/3/2</id>
            <config>
              <id>GigabitEthernet0/3/2</id>
              <interface>GigabitEthernet0/3/2</interface>
              <subinterface>0</subinterface>
            </config>
          </interface>
        </interfaces>
        <protocols>
          <protocol>
            <identifier>oc-pol-types:OSPF</identifier>
            <name>22</name>
            <config>
              <identifier>oc-pol-types:OSPF</identifier>
              <name>22</name>
            </config>
            <ospfv2>
              <global>
                <config>
                  <router-id>172.16.1.3</router-id>
                </config>
              </global>
            </ospfv2>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:STATIC</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:STATIC</identifier>
              <name>default</name>
            </config>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
              <name>default</name>
            </config>
          </protocol>
        </protocols>
      </network-instance>
    </network-instances>
  </data>
</rpc-reply>

'@
$g2 = @'
<edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:OSPF</identifier>
              <name>22</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:OSPF</identifier>
                <name>22</name>
              </config>
              <ospfv2>
                <areas>
                  <area>
                    <identifier>0.0.0.0</identifier>
                    <config>
                      <identifier>0.0.0.0</identifier>
                    </config>
                    <interfaces>
                      <interface>
                        <id>GigabitEthernet0/3/0</id>
                        <config>
                          <id>GigabitEthernet0/3/0</id>
                          <authentication-type>MD5</authentication-type>
                        </config>
                      </interface>
                    </interfaces>
                  </area>
                </areas>
              </ospfv2>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
</edit-config>
'@
$ws.Range("F2").Value = $f2
$ws.Range("G2").Value = $g2
